# UMG_Stammdaten_MAPPING.xlsx — fix UTF-8 encoding of FLAT paths and clean
# up stray manual entries on the "Mapping CSV2openEHR" sheet, then restore
# the per-row dropdown (list) data validation.
#
# Commit message: "Fixed Upload of OPT in correct UTF-8 encoding and added
# some error handling"

$wb = $excel.ActiveWorkbook

$mapSheet = $wb.Worksheets.Item("Mapping CSV2openEHR")
$pathSheet = $wb.Worksheets.Item("FLAT_Paths")

# --- 1) Correct the mangled German umlauts ("pr_fix" -> "präfix", etc.) on
#        the FLAT_Paths lookup sheet. The shared-string table is rebuilt on
#        save, so editing the cell text here is sufficient to fix every
#        place that string is used.
$pathSheet.Range("A5").Value2  = "bericht/context/umg_personenname<<index>>/name_strukturiert/präfix"
$pathSheet.Range("A6").Value2  = "bericht/context/umg_personenname<<index>>/name_strukturiert/art_des_präfix"
$pathSheet.Range("A18").Value2 = "bericht/context/umg_adresse<<index>>/straße"

# --- 2) Remove stray values that had been typed directly into the mapping
#        column/index cells instead of being picked from the validation
#        dropdown (rows 5 and 8).
$mapSheet.Range("B5:C5").ClearContents()
$mapSheet.Range("B8").ClearContents()

# --- 3) (Re-)apply the "choose a FLAT path" list validation for every
#        mapping row so invalid/free-typed entries are rejected again.
for ($r = 2; $r -le 17; $r++) {
    $cell = $mapSheet.Cells.Item($r, 2)
    $cell.Validation.Add(3, 1, 1, "=FLAT_Paths!`$A`$2:`$A`$93")
}

Write-Output "UMG_Stammdaten_MAPPING: encoding fixed, stray cells cleared, validation restored"
